$d = $word.ActiveDocument

# 1. Remove the "Meta description: ..." paragraph that follows the title.
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Delete()

# 2. Replace the final paragraph (the image-prompt paragraph) with two
#    paragraphs: a bold title line followed by the italic meta-description
#    text that used to live near the top of the document.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Absolootly Mad Mega Moolah for Free - Game Review</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Want to play Absolootly Mad Mega Moolah for free? Check our review of this jackpot game, featuring 4 progressive jackpots and up to 56 free spins.</w:t></w:r></w:p>'
[void]$lastPara.Range.InsertXML($xml)
